$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.915.76'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.551.02'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.38%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '206.34'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("E7").Value = '  -0.41%  '
$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '22.26'
$c.Style = 'Normal'
$ws.Range("E8").Value = '  +3.54%  '
$ws.Range("E9").Value = '  -0.43%  '
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.0587'
$c.Style = 'Normal'
$ws.Range("E10").Value = '  +0.63%  '
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.0855'
$c.Style = 'Normal'
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").Value = '1.770.60'
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D13").Value = '1.550.17'
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("E14").Value = '  +0.73%  '
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '0.518'
$c.Style = 'Normal'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '26.905.79'
$ws.Range("E16").Value = '  -0.11%  '
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '61.58'
$c.Style = 'Normal'
$ws.Range("E17").Value = '  -0.35%  '
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '217.20'
$c.Style = 'Normal'
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("E19").Value = '  +1.52%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("E21").Value = '  -0.39%  '
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '4.06'
$c.Style = 'Normal'
$ws.Range("E22").Value = '  +0.11%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '9.24'
$c.Style = 'Normal'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("E24").Value = '  -0.70%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '153.99'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("E26").Value = '  -0.82%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '14.93'
$c.Style = 'Normal'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("E28").Value = '  +0.82%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  -0.51%  '
$ws.Range("D33").Value = '1.418.82'
$ws.Range("E33").Value = '  +3.53%  '
$ws.Range("E34").Value = '  +3.59%  '
$ws.Range("E35").Value = '  +1.96%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  +0.38%  '
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '0.524'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("E42").Value = '  +3.48%  '
$ws.Range("E43").Value = '  +3.28%  '
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range("E44").Value = '  +1.54%  '
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '64.59'
$c.Style = 'Normal'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("E46").Value = '  +1.22%  '
$ws.Range("D47").Value = '1.684.96'
$ws.Range("E47").Value = '  -0.36%  '
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '87.20'
$c.Style = 'Normal'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("E49").Value = '  +3.97%  '
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '0.0515'
$c.Style = 'Normal'
$ws.Range("E50").Value = '  +1.68%  '
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '0.0959'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  +0.34%  '
